$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B4"  = 7.733999999999999
    "A9"  = -21.74
    "B9"  = 6.068000000000001
    "D9"  = -7.852000000000001
    "A18" = -21.995
    "A20" = -20.511
    "B23" = 7.62
    "B24" = 5.149000000000001
    "B26" = 5.223000000000001
    "A27" = -21.791
    "D32" = -7.226999999999999
    "B34" = 7.215999999999999
    "B35" = 8.031000000000001
    "D38" = -7.867
    "D45" = -7.503000000000002
    "B48" = 5.433999999999999
    "D51" = -8.409000000000001
    "B52" = 5.544
    "D57" = -8.145999999999999
    "D64" = -7.827
    "B66" = 5.172
    "B67" = 5.374
    "A69" = -21.521
    "A76" = -20.468
    "B80" = 8.606
    "A82" = -21.81
    "D93" = -7.246
    "B99" = 5.206999999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
